$wb = $excel.ActiveWorkbook

# Rename the first sheet from "nad" to "nadp"
$nadp = $wb.Worksheets.Item("nad")
$nadp.Name = "nadp"

# Add the new "dna" data (columns F:H, rows 2-25)
$dna = $wb.Worksheets.Item("dna")

$data = @(
    @(2, 45479, 47505, 47532),
    @(3, 594246, 592748, 607592),
    @(4, 1094433, 1084946, 1085169),
    @(5, 2116238, 2083453, 2167345),
    @(6, 3969176, 4062854, 4039748),
    @(7, 7683769, 7525415, 7405485),
    @(8, 14950687, 15095009, 15088696),
    @(9, 32258614, 31554882, 32605324),
    @(10, 6425032, 6455568, 6534228),
    @(11, 7833076, 5946024, 7100669),
    @(12, 6442415, 6468564, 6297042),
    @(13, 5860193, 6294057, 5660226),
    @(14, 6428410, 6484686, 6529611),
    @(15, 5046831, 5823830, 5855774),
    @(16, 3763118, 2974154, 2764465),
    @(17, 3106065, 3690916, 4111486),
    @(18, 2809087, 3082595, 3912677),
    @(19, 4546006, 3136722, 4037084),
    @(20, 5163304, 4696128, 5921674),
    @(21, 2356400, 2458309, 3284193),
    @(22, 2649190, 3088234, 3904475),
    @(23, 2375775, 2395540, 3545720),
    @(24, 3221743, 2444418, 2322598),
    @(25, 2389760, 2197585, 3404384)
)

foreach ($row in $data) {
    $r = $row[0]
    $dna.Cells.Item($r, 6).Value = $row[1]
    $dna.Cells.Item($r, 7).Value = $row[2]
    $dna.Cells.Item($r, 8).Value = $row[3]
}

# Update the selection / active cell on each sheet to match the recorded view state
$dna.Range("F18:H25").Select() | Out-Null
$nadp.Range("C31").Select() | Out-Null
